# Update QT counts to HEWR only.
# - Adds a new "Area" column (S) identifying the survey/estimate area for
#   each sourced row (HEWR vs. HEWR+rugged west).
# - Revises the 2013 (row 14) and 2016 (row 17) MinCount_ADULTMF /
#   MinCount_CALFMF figures (MinCount total recalculates via its SUM formula).
# - Expands the note on the 2013 estimate with additional context.
# - Leaves the active selection on G17, matching the author's last position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Area" header/column.
$ws.Range("S1").Value = "Area"
$ws.Range("S3").Value = "HEWR"
$ws.Range("S9").Value = "HEWR"
$ws.Range("S14").Value = "HEWR+rugged west"

# Row 14 (2013): MinCount_ADULTMF / MinCount_CALFMF updated; MinCount (C14)
# recalculates automatically from =SUM(D14:E14).
$ws.Range("D14").Value = 105
$ws.Range("E14").Value = 9

# Row 17 (2016): same kind of update.
$ws.Range("D17").Value = 32
$ws.Range("E17").Value = 7

# Expanded comment for the 2013 estimate.
$ws.Range("Q14").Value = 'Estimate was 129 in report. Note (BM): We used the population estimate based strictly on the caribou seen within the survey area. // MC: This is the estimate for the entire population when we correct for 24% of the population that may have been outside the census area (i.e. four of the 17 radio-collared caribou were outside the census area.). Note (BM): Another MC population estimate is given; 129 caribou, with a lower and upper 95% CI of (114, 129) and a SCF of 0.76. "This is the estimate for the entire population when we correct for 24% of the population that may have been outside the census area (i.e. four of the 17 radio-collared caribou were outside the census area.)" // Note (BM): Confidence level for population estimate is 95%.'

# Match the author's final selection.
$ws.Range("G17").Select()
